$d = $word.ActiveDocument

# Locate the run containing the Russian text "Моя" that needs to become "My"
$findRange = $d.Content
$findRange.Find.Execute("Моя", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$para = $findRange.Paragraphs(1)

# Remove the old text, leaving just the paragraph mark behind. With the
# paragraph now empty, setting the language on its range stamps the
# paragraph-mark run properties (w:pPr/w:rPr/w:lang).
$findRange.Text = ""
$para.Range.LanguageID = "en-US"

# Insert the new English text at the start of the (now empty) paragraph.
$insertStart = $para.Range.Start
$insertRange = $d.Range($insertStart, $insertStart)
$insertRange.Text = "My"

# Stamp the language on the freshly inserted run too (w:r/w:rPr/w:lang).
$runRange = $d.Range($insertStart, $insertStart + 2)
$runRange.LanguageID = "en-US"
